# Update cryptocurrency price (D) and 1h volume change (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.832.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.09'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.67%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("E8").Value = '  -0.55%  '

$ws.Range("E9").Value = '  -0.94%  '

$ws.Range("E10").Value = '  +0.61%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("E12").Value = '  +0.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.627.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '

$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.836.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.50%  '

$ws.Range("E20").Value = '  -0.15%  '

$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("E24").Value = '  -1.54%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.07%  '

$ws.Range("E26").Value = '  -0.21%  '

$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("E29").Value = '  +0.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0508'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.16%  '

$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("E33").Value = '  -0.31%  '

$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.274.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.39%  '

$ws.Range("E36").Value = '  +0.15%  '

$ws.Range("E37").Value = '  -1.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.529'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.07%  '

$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.803'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.780.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.64%  '

$ws.Range("E44").Value = '  -5.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = '  -1.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0968'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.51%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.65%  '

$ws.Range("E51").Value = '  -0.07%  '
